# Weekly update: insert two new price records (most recent week, 2023-09-25)
# at the top of the data block, pushing all existing rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current first data row (row 697),
# shifting every existing record down by two rows.
$ws.Rows.Item(697).Insert()
$ws.Rows.Item(697).Insert()

# --- New row 697 ---
$ws.Cells.Item(697,1).Value = 10
$ws.Cells.Item(697,2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(697,3).Value = 'La Araucanía'
$ws.Cells.Item(697,4).Value = 45194
$ws.Cells.Item(697,5).Value = 9
$ws.Cells.Item(697,6).Value = 100112023
$ws.Cells.Item(697,7).Value = 'Brócoli'
$ws.Cells.Item(697,8).Value = 'Sin especificar'
$ws.Cells.Item(697,9).Value = 'Primera'
$ws.Cells.Item(697,10).Value = 1400
$ws.Cells.Item(697,11).Value = 1200
$ws.Cells.Item(697,12).Value = 1200
$ws.Cells.Item(697,13).Value = 1200
$ws.Cells.Item(697,14).Value = '$/unidad'
$ws.Cells.Item(697,15).Value = 'Provincia del Elquí'
$ws.Cells.Item(697,16).Value = 1200
$ws.Cells.Item(697,17).Value = 1
$ws.Cells.Item(697,18).Value = 'Hortaliza'

# --- New row 698 ---
$ws.Cells.Item(698,1).Value = 10
$ws.Cells.Item(698,2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(698,3).Value = 'La Araucanía'
$ws.Cells.Item(698,4).Value = 45194
$ws.Cells.Item(698,5).Value = 9
$ws.Cells.Item(698,6).Value = 100112023
$ws.Cells.Item(698,7).Value = 'Brócoli'
$ws.Cells.Item(698,8).Value = 'Sin especificar'
$ws.Cells.Item(698,9).Value = 'Primera'
$ws.Cells.Item(698,10).Value = 2500
$ws.Cells.Item(698,11).Value = 1200
$ws.Cells.Item(698,12).Value = 1300
$ws.Cells.Item(698,13).Value = 1260
$ws.Cells.Item(698,14).Value = '$/unidad'
$ws.Cells.Item(698,15).Value = 'Región Metropolitana'
$ws.Cells.Item(698,16).Value = 1260
$ws.Cells.Item(698,17).Value = 1
$ws.Cells.Item(698,18).Value = 'Hortaliza'
